$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Age Group")
$ws.Range("B2").Value = 44802
$ws.Range("C2").Value = 3840
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 7.93
$ws.Range("F2").Value = 7.44
$ws.Range("G2").Value = 0.08
$ws.Range("B3").Value = 85970
$ws.Range("C3").Value = 8609
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 15.23
$ws.Range("F3").Value = 16.68
$ws.Range("G3").Value = 0.31
$ws.Range("B4").Value = 87684
$ws.Range("C4").Value = 8244
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 15.53
$ws.Range("F4").Value = 15.97
$ws.Range("G4").Value = 0.7
$ws.Range("B5").Value = 85290
$ws.Range("C5").Value = 8536
$ws.Range("D5").Value = 46
$ws.Range("E5").Value = 15.11
$ws.Range("F5").Value = 16.54
$ws.Range("G5").Value = 1.79
$ws.Range("B6").Value = 92139
$ws.Range("C6").Value = 7998
$ws.Range("D6").Value = 131
$ws.Range("E6").Value = 16.32
$ws.Range("F6").Value = 15.5
$ws.Range("G6").Value = 5.1
$ws.Range("B7").Value = 81088
$ws.Range("C7").Value = 6151
$ws.Range("D7").Value = 407
$ws.Range("E7").Value = 14.36
$ws.Range("F7").Value = 11.92
$ws.Range("G7").Value = 15.86
$ws.Range("B8").Value = 50947
$ws.Range("C8").Value = 3880
$ws.Range("D8").Value = 627
$ws.Range("E8").Value = 9.02
$ws.Range("F8").Value = 7.52
$ws.Range("G8").Value = 24.43
$ws.Range("B9").Value = 33456
$ws.Range("C9").Value = 4322
$ws.Range("D9").Value = 1328
$ws.Range("E9").Value = 5.93
$ws.Range("F9").Value = 8.37
$ws.Range("G9").Value = 51.73
$ws.Range("B10").Value = 3271
$ws.Range("C10").Value = 32
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.58
$ws.Range("F10").Value = 0.06
$ws.Range("G10").Value = 0

$ws = $wb.Worksheets.Item("Gender")
$ws.Range("B2").Value = 324072
$ws.Range("C2").Value = 26952
$ws.Range("D2").Value = 1252
$ws.Range("E2").Value = 57.39
$ws.Range("F2").Value = 52.22
$ws.Range("G2").Value = 48.77
$ws.Range("B3").Value = 229032
$ws.Range("C3").Value = 23749
$ws.Range("D3").Value = 1261
$ws.Range("E3").Value = 40.56
$ws.Range("F3").Value = 46.01
$ws.Range("G3").Value = 49.12
$ws.Range("B4").Value = 11543
$ws.Range("C4").Value = 911
$ws.Range("D4").Value = 54
$ws.Range("E4").Value = 2.04
$ws.Range("F4").Value = 1.77
$ws.Range("G4").Value = 2.1

$ws = $wb.Worksheets.Item("Race")
$ws.Range("B2").Value = 4767
$ws.Range("C2").Value = 738
$ws.Range("D2").Value = 13
$ws.Range("E2").Value = 0.84
$ws.Range("F2").Value = 1.43
$ws.Range("G2").Value = 0.51
$ws.Range("B3").Value = 49068
$ws.Range("C3").Value = 6099
$ws.Range("D3").Value = 370
$ws.Range("E3").Value = 8.69
$ws.Range("F3").Value = 11.82
$ws.Range("G3").Value = 14.41
$ws.Range("B4").Value = 76323
$ws.Range("C4").Value = 8989
$ws.Range("D4").Value = 337
$ws.Range("E4").Value = 13.52
$ws.Range("F4").Value = 17.42
$ws.Range("G4").Value = 13.13
$ws.Range("B5").Value = 110803
$ws.Range("C5").Value = 12292
$ws.Range("D5").Value = 185
$ws.Range("E5").Value = 19.62
$ws.Range("F5").Value = 23.82
$ws.Range("G5").Value = 7.21
$ws.Range("B6").Value = 323686
$ws.Range("C6").Value = 23494
$ws.Range("D6").Value = 1662
$ws.Range("E6").Value = 57.33
$ws.Range("F6").Value = 45.52
$ws.Range("G6").Value = 64.74

$ws = $wb.Worksheets.Item("Ethnicity")
$ws.Range("B2").Value = 18617
$ws.Range("C2").Value = 5805
$ws.Range("D2").Value = 51
$ws.Range("E2").Value = 3.3
$ws.Range("F2").Value = 11.25
$ws.Range("G2").Value = 1.99
$ws.Range("B3").Value = 215148
$ws.Range("C3").Value = 17255
$ws.Range("D3").Value = 1357
$ws.Range("E3").Value = 38.1
$ws.Range("F3").Value = 33.43
$ws.Range("G3").Value = 52.86
$ws.Range("B4").Value = 330882
$ws.Range("C4").Value = 28552
$ws.Range("D4").Value = 1159
$ws.Range("E4").Value = 58.6
$ws.Range("F4").Value = 55.32
$ws.Range("G4").Value = 45.15
